$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the three tasks that were dropped entirely (delete from the
#    bottom up so earlier row numbers stay valid).
#       Row 16: "Hyper parameter tuning"
#       Row 12: "Translation (follow up with Sumit)"
#       Row 8 : "Data Augmentation"
#    The table (ListObject) auto-shrinks with the sheet rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# 2. Update remaining rows with the new Status / Column1 / Column2 content.
# ---------------------------------------------------------------------------

# Row 6 - Translation, clean garbage... : drop the "Doubt" Column1 flag
$ws.Cells.Item(6, 5).ClearContents()

# Row 8 (was row 9) - Data Augmentation and Development of New Data Frames
$ws.Cells.Item(8, 4).Value2 = "Almost Complete"
$ws.Cells.Item(8, 6).Value2 = "Need to graph before and after and provide examples of augmentation and translation"

# Row 10 (was row 11) - Doc2Vec Exploration with Gensim
$ws.Cells.Item(10, 4).Value2 = "Pending"

# Row 12 (was row 14) - Feature Engineering...
$ws.Cells.Item(12, 3).Value2 = "Abhik"
$ws.Cells.Item(12, 4).Value2 = "Completed"
$ws.Cells.Item(12, 5).ClearContents()

# Row 13 (was row 15) - Machine learning models (SVM, Logistic, XGBoost)
$ws.Cells.Item(13, 3).Value2 = "Sayantika"
$ws.Cells.Item(13, 4).Value2 = "Almost Complete"
$ws.Cells.Item(13, 5).ClearContents()
$ws.Cells.Item(13, 6).Value2 = "Sayantika to update predictions, examples and graphing"

# Row 14 (was row 17) - Deep Learning
$ws.Cells.Item(14, 4).Value2 = "Almost Complete"
$ws.Cells.Item(14, 6).Value2 = "Need to append graphs and prediction examples"

# Row 15 (was row 18) - Create Report
$ws.Cells.Item(15, 4).Value2 = "Pending"
$ws.Cells.Item(15, 6).Value2 = "We need to connect on Tuesday and put an outline together"

# Row 16 (was row 19) - Validate the deliverables against Requirements
$ws.Cells.Item(16, 4).Value2 = "Pending"

# Row 17 (was row 20) - Rule based system
$ws.Cells.Item(17, 3).Value2 = "Abhijit"
$ws.Cells.Item(17, 4).Value2 = "Pending"

# Row 18 (was row 21) - Improvement of visualizations
$ws.Cells.Item(18, 3).Value2 = "Arun"
$ws.Cells.Item(18, 4).Value2 = "Pending"

# Row 19 (was row 22) - Move Translation to Goslate -> rename + new data
$ws.Cells.Item(19, 2).Value2 = "Move Translation to Goslate/Improving translations"
$ws.Cells.Item(19, 3).Value2 = "Abhik"
$ws.Cells.Item(19, 4).Value2 = "Almost Complete"
$ws.Cells.Item(19, 6).Value2 = "Good to have. Not Mandatory"

# ---------------------------------------------------------------------------
# 3. Append two brand-new rows (20 & 21) for WEEK3.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value2 = "WEEK3"
$ws.Cells.Item(20, 2).Value2 = "Use tensor board for visualizations and find hyperparameters"
$ws.Cells.Item(20, 3).Value2 = "Abhik"
$ws.Cells.Item(20, 4).Value2 = "Pending"

$ws.Cells.Item(21, 1).Value2 = "WEEK3"
$ws.Cells.Item(21, 2).Value2 = "Stratified Training data split"
$ws.Cells.Item(21, 3).Value2 = "Abhijit"
$ws.Cells.Item(21, 4).Value2 = "Pending"

# Match the "Week" column centred style used by the rest of column A
# (style index 1 = horizontal-center alignment).
$xlCenter = -4108
$ws.Range("A20:A21").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# 4. Resize the table to the new extent and refresh the selection.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F21"))

$ws.Range("B24").Select()
